# Apply cryptos list price/volume updates (and a few coin-row swaps)
# Generated from the canonical OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.797.52"
$ws.Range("E2").Value = "  +5.01%  "
$ws.Range("D3").Value = "3.263.92"
$ws.Range("E3").Value = "  +5.35%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'579.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.09%  "
$ws.Range("D6").Value = "'181.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.41%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").Value = "3.263.49"
$ws.Range("E9").Value = "  +5.43%  "
$ws.Range("D10").Value = "'0.133"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.28%  "
$ws.Range("D11").Value = "'6.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.11%  "
$ws.Range("D12").Value = "'0.418"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.21%  "
$ws.Range("D13").Value = "3.835.40"
$ws.Range("E13").Value = "  +5.56%  "
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "'28.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.91%  "
$ws.Range("D16").Value = "67.782.76"
$ws.Range("E16").Value = "  +5.16%  "
$ws.Range("E17").Value = "  +6.15%  "
$ws.Range("D18").Value = "3.266.15"
$ws.Range("E18").Value = "  +5.07%  "
$ws.Range("D19").Value = "'5.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.13%  "
$ws.Range("D20").Value = "'13.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.31%  "
$ws.Range("D21").Value = "'375.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.76%  "
$ws.Range("D22").Value = "'7.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.22%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'71.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.23%  "
$ws.Range("D25").Value = "'0.511"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.93%  "
$ws.Range("E26").Value = "  +9.33%  "
$ws.Range("D27").Value = "'9.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.11%  "
$ws.Range("E28").Value = "  +4.31%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +5.16%  "
$ws.Range("D31").Value = "'5.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.50%  "
$ws.Range("E32").Value = "  +6.67%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.37%  "
$ws.Range("D35").Value = "'6.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.74%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'163.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.91%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.18%  "
$ws.Range("D38").Value = "'0.850"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.74%  "
$ws.Range("D39").Value = "'1.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.07%  "
$ws.Range("D40").Value = "'6.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.26%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'26.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.36%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'4.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.56%  "
$ws.Range("D43").Value = "'2.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.68%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.701.01"
$ws.Range("E44").Value = "  +3.49%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'351.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.81%  "
$ws.Range("D46").Value = "'25.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.77%  "
$ws.Range("D47").Value = "'40.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.42%  "
$ws.Range("D48").Value = "'0.0680"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.96%  "
$ws.Range("D49").Value = "'0.0281"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.92%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.19%  "
$ws.Range("D51").Value = "'0.102"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.96%  "
